$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s-val rows (B:E recomputed; G = running sum B+C+D+E) after filtering save games.
$newVals = @{
    2 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    5 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    7 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    8 = @(0.003994804209775715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0.9710897032086083)
    9 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    10 = @(0.127881588408715, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.4499806748245367)
    11 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.7429408310145853)
    12 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    13 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    14 = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.0949764874075476)
    15 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    16 = @(1.459612070389937, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 2.4267069693887695)
    17 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217081)
    18 = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.271902929317955)
    19 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    20 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    21 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    22 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    23 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    24 = @(0.3048080303191223, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 1.0015170202094372)
    25 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    26 = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 2.290389397800092)
    27 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 31.612965916961354)
    28 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    29 = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.271902929317955)
    30 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    31 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    32 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.7403346288415715)
    33 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    34 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.7817111568057586)
    35 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    36 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    37 = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.9198672729249926)
    38 = @(0.127881588408715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.837881874639075)
    39 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    40 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.7817111568057586)
    41 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    42 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    43 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    44 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    45 = @(1.459612070389937, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 2.4267069693887695)
    46 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    47 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 17.459443432731913)
    48 = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 2.290389397800092)
    49 = @(0.04763786555579896, 0.002777888934908601, 3.900430680208489, 0.496779210170732, 4.4476256448699285)
    50 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217081)
}

foreach ($row in $newVals.Keys) {
    $vals = $newVals[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E: IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G: sum
}

